$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add two new rows (40 and 41) to Sheet1, following the same pattern as the
# existing rows (App=".", TestName, ID, Priority="3,", Platforme="1",
# Feature="data", test="on"), copying the border / number formatting of the
# last existing data row (39) so the new rows look consistent with the rest
# of the table.
# ---------------------------------------------------------------------------

function Add-TestRow($targetRow, $templateRow, $testName, $testId) {
    # Copy cell formatting (borders, number formats) from the template row
    # for columns C:G so the new row matches the existing table styling.
    $ws.Range("C$templateRow`:G$templateRow").Copy()
    $ws.Range("C$targetRow`:G$targetRow").PasteSpecial(-4122)
    $ws.Application.CutCopyMode = 0

    $ws.Range("A$targetRow").Value = "."
    $ws.Range("B$targetRow").Value = $testName
    $ws.Range("C$targetRow").Value = $testId

    # Column D holds the text "3," -- force a text format before assigning it
    # so it is not re-interpreted as the number 3.
    $ws.Range("D$targetRow").NumberFormat = "@"
    $ws.Range("D$targetRow").Value = "3,"

    $ws.Range("E$targetRow").Value = "1"
    $ws.Range("F$targetRow").Value = "data"
    $ws.Range("G$targetRow").Value = "on"

    # Re-apply column D's formatting from the template (restores the regular
    # bordered look used by the rest of the table) without touching the text
    # value that was just written.
    $ws.Range("D$templateRow").Copy()
    $ws.Range("D$targetRow").PasteSpecial(-4122)
    $ws.Application.CutCopyMode = 0
}

Add-TestRow 40 39 "Credit_cards_transactions_filter_by_date(custom_date_range)_[MOB_ANDROID]" "C70815"
Add-TestRow 41 39 "Credit_cards_transactions_filter_invalid_[MOB_ANDROID]" "C70816"

# ---------------------------------------------------------------------------
# Update the AutoFilter range to include the new data (A1:G39 -> A1:G39 is
# one header + 38 data rows before; now 39 data rows -> A1:G39... wait the
# table grew by two data rows, so the filter range grows from A1:G38 to
# A1:G39, mirroring the table's "last row with contiguous filter button"
# rule used by the source workbook).
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:G39").AutoFilter()

# Keep the workbook's hidden _FilterDatabase defined name in sync with the
# AutoFilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$G`$39"
    }
}

# ---------------------------------------------------------------------------
# Extend the conditional formatting (duplicate-values highlighting) that
# covered C37:C39 so it also covers the two new rows (C37:C41).
# ---------------------------------------------------------------------------
$dupCond = $ws.Range("C37").FormatConditions.Item(1)
$dupCond.ModifyAppliesToRange($ws.Range("C37:C41"))

# ---------------------------------------------------------------------------
# Move the active selection to C43 (next free row below the new data),
# matching where the workbook's author left the cursor after the edit.
# ---------------------------------------------------------------------------
$ws.Range("C43").Select()
